$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Swap the data of row 8 (Haywood Highsmith) and row 9 (Kyle Lowry) so
#    that Kyle Lowry's record now appears on row 8 and Haywood Highsmith's
#    record appears on row 9 (columns B..K). Column A, the roster index,
#    stays untouched.
#
#    We use Range.Copy through a scratch row far away from the real data so
#    that the original cell types (numbers vs. text, e.g. the "Exp" column
#    holds numeric-looking text like "2"/"16" that must not be turned into
#    real numbers) and formatting survive the swap exactly, and then wipe
#    the scratch row afterwards so no stray data / dimension growth remains.
# ---------------------------------------------------------------------------

$scratch = 1000

$ws.Range("B8:K8").Copy($ws.Range("B" + $scratch + ":K" + $scratch))
$ws.Range("B9:K9").Copy($ws.Range("B8:K8"))
$ws.Range("B" + $scratch + ":K" + $scratch).Copy($ws.Range("B9:K9"))
$ws.Range("B" + $scratch + ":K" + $scratch).Clear()

# ---------------------------------------------------------------------------
# 2) Row 18 (Omer Yurtseven) was missing a jersey number; add it.
# ---------------------------------------------------------------------------
$ws.Range("B18").Value = 77
